$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - field types
$ws.Range("A1").Value = "int(11)"
$ws.Range("B1").Value = "varchar(50)"
$ws.Range("C1").Value = "varchar(50)"
$ws.Range("D1").Value = "varchar(50)"
$ws.Range("E1").Value = "varchar(50)"
$ws.Range("F1").Value = "varchar(50)"
$ws.Range("G1").Value = "varchar(50)"
$ws.Range("H1").Value = "varchar(6)"
$ws.Range("I1").Value = "int(255)"
$ws.Range("J1").Value = "int(255)"

# Row 2 - field labels
$ws.Range("A2").Value = "ID (입력x)"
$ws.Range("B2").Value = "공장코드"
$ws.Range("C2").Value = "법인코드"
$ws.Range("D2").Value = "계정코드"
$ws.Range("E2").Value = "version코드"
$ws.Range("F2").Value = "사업장코드"
$ws.Range("G2").Value = "bom코드"
$ws.Range("H2").Value = "년월 ex) 200001"
$ws.Range("I2").Value = "수량"
$ws.Range("J2").Value = "금액"
